$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (U, V): sundown program additions ----------------------
$ws.Range("U1").Value = "sunDownDate"
$ws.Range("V1").Value = "sunDownLength"

# --- Full factorial grid: NUM_INSTRUCTOR (C) x s_o_c (D) ------------------------
# C cycles fastest through 10,30,50,70,90 ; D steps through 25,50,75 -> rows 2..16
$cVals = @(10, 30, 50, 70, 90)
$dVals = @(25, 50, 75)

$r = 2
foreach ($d in $dVals) {
    foreach ($c in $cVals) {
        $ws.Cells.Item($r, 1).Value = 50      # A NUM_AIRCRAFT
        $ws.Cells.Item($r, 2).Value = 50      # B NUM_STUDENT
        $ws.Cells.Item($r, 3).Value = $c      # C NUM_INSTRUCTOR
        $ws.Cells.Item($r, 4).Value = $d      # D s_o_c
        $ws.Cells.Item($r, 5).Value = 42      # E rl
        $ws.Cells.Item($r, 6).Value = 720     # F ip
        $ws.Cells.Item($r, 7).Value = 0.035   # G attrit
        $ws.Cells.Item($r, 9).Value = 7000    # I sleplimit
        $ws.Cells.Item($r, 10).Value = 720    # J et_af
        $ws.Cells.Item($r, 11).Value = 240    # K et_av
        $ws.Cells.Item($r, 12).Value = 360    # L et_puls
        $ws.Cells.Item($r, 13).Value = 720    # M rt_af
        $ws.Cells.Item($r, 14).Value = 480    # N rt_av
        $ws.Cells.Item($r, 15).Value = 240    # O rt_puls
        $ws.Cells.Item($r, 16).Value = $true  # P SLEP_or_not
        $ws.Cells.Item($r, 17).Value = 8      # Q SLEPspots
        $ws.Cells.Item($r, 18).Value = $true  # R Stagger
        $ws.Cells.Item($r, 19).Value = 14400  # S addHours
        $r = $r + 1
    }
}

# time_line (H) is the same formula for every row -> fill as one range so Excel
# collapses it into a single shared formula (H2 master, ref H2:H16).
$ws.Range("H2:H16").Formula = "=24*365*50"

# TTR (T): rows 2-3 keep their own individual formulas (unchanged from before),
# rows 4-16 were filled together and become a shared-formula group (T4 master,
# ref T4:T16).
$ws.Range("T2").Formula = "=24*30*9"
$ws.Range("T3").Formula = "=24*30*9"
$ws.Range("T4:T16").Formula = "=24*30*9"
